# Updated cryptos list on Tue Nov 26 06:43:35 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row, and
# re-sorts a handful of near-tied rows (USDe/Kaspa, Filecoin/MantraDAO/
# VeChain/Stacks/dogwifhat) by swapping their Coin/Link/Price/Volume cells.
#
# Cells whose new text looks like a plain number (e.g. "644.09") are written
# with a leading apostrophe so Excel stores them as text instead of
# re-interpreting them as numeric values (matching the original inlineStr
# cells, which are all plain text). The apostrophe flips the cell's style to
# a "quote prefix" style, so we immediately reset those cells back to the
# "Normal" style afterwards to match the unstyled cells elsewhere in the
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textUpdates = @{
    "D2"  = "94.782.17"
    "E2"  = "  -3.54%  "
    "D3"  = "3.440.27"
    "E3"  = "  +1.63%  "
    "E4"  = "  -0.01%  "
    "E5"  = "  -5.32%  "
    "E7"  = "  -1.26%  "
    "E8"  = "  -3.84%  "
    "E9"  = "  +0.10%  "
    "E10" = "  -5.93%  "
    "D11" = "3.438.77"
    "E11" = "  +1.65%  "
    "E12" = "  -4.05%  "
    "E13" = "  -1.75%  "
    "E14" = "  +2.78%  "
    "D15" = "94.687.34"
    "E15" = "  -3.21%  "
    "D16" = "4.081.35"
    "E16" = "  +1.72%  "
    "E17" = "  -0.97%  "
    "E18" = "  -9.47%  "
    "D19" = "3.437.06"
    "E19" = "  +1.18%  "
    "E20" = "  -2.37%  "
    "E21" = "  +6.26%  "
    "E22" = "  -3.45%  "
    "E23" = "  -1.45%  "
    "E24" = "  -4.56%  "
    "E25" = "  -3.59%  "
    "E26" = "  -4.80%  "
    "E27" = "  -2.25%  "
    "D29" = "3.623.47"
    "E29" = "  +1.54%  "
    "E30" = "  +0.27%  "
    "E31" = "  +0.13%  "
    "E32" = "  +6.18%  "
    "E33" = "  -3.32%  "
    "E34" = "  -0.40%  "
    "E35" = "  -4.63%  "
    "E36" = "  +3.70%  "
    "E37" = "  -0.78%  "
    "E38" = "  +4.07%  "
    "E39" = "  -2.36%  "
    "E40" = "  -0.75%  "
    "B41" = "Kaspa"
    "C41" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "E41" = "  -0.43%  "
    "B42" = "USDe"
    "C42" = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
    "E42" = "  +0.02%  "
    "E43" = "  +8.58%  "
    "E44" = "  -1.29%  "
    "E45" = "  -0.40%  "
    "B46" = "MantraDAO"
    "C46" = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
    "E46" = "  +0.16%  "
    "B47" = "Filecoin"
    "C47" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
    "E47" = "  +2.52%  "
    "B48" = "VeChain"
    "C48" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "E48" = "  -3.36%  "
    "B49" = "Stacks"
    "C49" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "E49" = "  -0.73%  "
    "B50" = "dogwifhat"
    "C50" = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
    "E50" = "  +3.84%  "
    "E51" = "  -0.99%  "
}

# Values that parse as a plain number and therefore need to be forced to
# text (leading apostrophe) so Excel keeps them as strings like the
# original cells.
$numericTextUpdates = @{
    "D5"  = "239.13"
    "D6"  = "644.09"
    "D10" = "0.981"
    "D13" = "41.95"
    "D14" = "6.29"
    "D18" = "8.39"
    "D20" = "17.62"
    "D21" = "11.69"
    "D22" = "0.498"
    "D23" = "503.17"
    "D24" = "3.26"
    "D26" = "6.57"
    "D27" = "94.55"
    "D28" = "12.01"
    "D30" = "11.76"
    "D33" = "0.138"
    "D34" = "0.999"
    "D35" = "0.180"
    "D36" = "29.91"
    "D37" = "0.556"
    "D38" = "552.54"
    "D39" = "7.74"
    "D41" = "0.152"
    "D42" = "1.00"
    "D43" = "0.922"
    "D45" = "1.73"
    "D46" = "3.69"
    "D47" = "5.70"
    "D48" = "0.0412"
    "D49" = "2.24"
    "D50" = "3.34"
    "D51" = "55.29"
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

foreach ($ref in $numericTextUpdates.Keys) {
    $ws.Range($ref).Value = "'" + $numericTextUpdates[$ref]
}

# Restore the default ("Normal") style on the cells we force-quoted above so
# they don't keep a "quote prefix" style that the original cells never had.
foreach ($ref in $numericTextUpdates.Keys) {
    $ws.Range($ref).Style = "Normal"
}
